$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Date Donated" column (H) for rows 8-12 with new donation dates.
# Order matches how the new shared strings were appended in the target workbook.
$ws.Range("H9").Value = "03/04/2022"
$ws.Range("H8").Value = "01/04/2022"
$ws.Range("H10").Value = "02/04/2022"
$ws.Range("H11").Value = "04/04/2022"
$ws.Range("H12").Value = "01/04/2021"

# Update the sheet view: scroll so column C is the left-most visible column,
# and move the active selection to H7.
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
$ws.Range("H7").Select()
